$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.549.69"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.882.53"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  +2.67%  "
$ws.Range("D6").Value = "'242.38"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07935"
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "'0.3120"
$ws.Range("E9").Value = "  +3.51%  "
$ws.Range("D10").Value = "'25.34"
$ws.Range("E10").Value = "  +6.58%  "
$ws.Range("D11").Value = "'0.08278"
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("D12").Value = "'0.7315"
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("D13").Value = "1.871.95"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "'5.291"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").Value = "'91.46"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").Value = "29.536.10"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "'5.950"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "'246.40"
$ws.Range("E18").Value = "  +4.63%  "
$ws.Range("D19").Value = "'0.000007897"
$ws.Range("E19").Value = "  +1.10%  "
$ws.Range("D20").Value = "'13.40"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").Value = "2.135.81"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'7.973"
$ws.Range("E23").Value = "  +6.25%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'0.1621"
$ws.Range("E25").Value = "  +14.27%  "
$ws.Range("D26").Value = "'163.38"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'9.079"
$ws.Range("E27").Value = "  +2.52%  "
$ws.Range("D28").Value = "'18.37"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("D30").Value = "'1.504"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "'4.399"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("D32").Value = "'4.117"
$ws.Range("E32").Value = "  +2.66%  "
$ws.Range("D33").Value = "'0.05294"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("D34").Value = "'1.963"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "'1.204"
$ws.Range("E35").Value = "  +3.64%  "
$ws.Range("D36").Value = "'0.7287"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'0.01874"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").Value = "1.233.75"
$ws.Range("E39").Value = "  +7.03%  "
$ws.Range("D40").Value = "'2.714"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("D41").Value = "'0.9118"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.228"
$ws.Range("E42").Value = "  +4.30%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'74.07"
$ws.Range("E43").Value = "  +5.87%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'102.02"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "2.031.78"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'1.796"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("D49").Value = "'2.938"
$ws.Range("E49").Value = "  +11.01%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'9.350"
$ws.Range("E51").Value = "  +2.62%  "
